$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.152.58'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '3.928.21'
$ws.Range('E3').Value = '  +3.69%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '470.73'
$ws.Range('E5').Value = '  +8.84%  '
$ws.Range('D6').Value = '145.88'
$ws.Range('E6').Value = '  +3.92%  '
$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +0.87%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.734'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E10').Value = '  +8.26%  '
$ws.Range('E11').Value = '  +7.90%  '
$ws.Range('D12').Value = '43.44'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').Value = '4.558.91'
$ws.Range('E13').Value = '  +3.95%  '
$ws.Range('D14').Value = '10.44'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = '15.16'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').Value = '3.912.11'
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').Value = '19.86'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('E19').Value = '  +1.82%  '
$ws.Range('D20').Value = '67.454.66'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').Value = '438.64'
$ws.Range('E21').Value = '  +6.97%  '
$ws.Range('D22').Value = '3.38'
$ws.Range('E22').Value = '  +3.37%  '
$ws.Range('E23').Value = '  -1.39%  '
$ws.Range('D24').Value = '87.91'
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('E25').Value = '  +7.26%  '
$ws.Range('D26').Value = '38.84'
$ws.Range('E26').Value = '  +5.24%  '
$ws.Range('D27').Value = '10.34'
$ws.Range('E27').Value = '  +5.63%  '
$ws.Range('D28').Value = '9.72'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').Value = '721.36'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').Value = '13.56'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '2.83'
$ws.Range('E32').Value = '  +2.81%  '
$ws.Range('D33').Value = '42.86'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '57.95'
$ws.Range('E34').Value = '  +3.34%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.151'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('D36').Value = '0.0₃0797'
$ws.Range('E36').Value = '  +17.74%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').Value = '5.36'
$ws.Range('E38').Value = '  -5.03%  '
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').Value = '3.05'
$ws.Range('E40').Value = '  +4.04%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '3.55'
$ws.Range('E41').Value = '  +6.49%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.142'
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '0.337'
$ws.Range('E43').Value = '  +3.26%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '2.57'
$ws.Range('E44').Value = '  -5.88%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '2.20'
$ws.Range('E46').Value = '  +5.66%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  +4.29%  '
$ws.Range('D48').Value = '147.45'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('D49').Value = '3.15'
$ws.Range('E49').Value = '  -5.64%  '
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').Value = '25.98'
$ws.Range('E51').Value = '  +3.96%  '
